# Insert a new row at position 213, shifting existing rows 213..265 down to 214..266,
# then populate the newly inserted row 213 with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(213).Insert()

$ws.Range("A213").Value = 8
$ws.Range("B213").Value = "Terminal La Palmera de La Serena"
$ws.Range("C213").Value = "Coquimbo"
$ws.Range("D213").Value = 44508
$ws.Range("E213").Value = 4
$ws.Range("F213").Value = 100114001
$ws.Range("G213").Value = "Papa"
$ws.Range("H213").Value = "Cardinal"
$ws.Range("I213").Value = "1a (cosecha)"
$ws.Range("J213").Value = 2400
$ws.Range("K213").Value = 11500
$ws.Range("L213").Value = 12000
$ws.Range("M213").Value = 11750
$ws.Range("N213").Value = "$/saco 25 kilos"
$ws.Range("O213").Value = "Provincia del Elquí"
$ws.Range("P213").Value = 470
$ws.Range("Q213").Value = 25
$ws.Range("R213").Value = "Hortaliza"
